# Generate Report for Handoff
#
# Regenerates the handoff XLIFF files for the four "Ready for handoff" source
# files (192f98ce..., 25d69de0..., 43854d5b..., 824acc8a...) on both the
# zh-cn and de-de worksheets: bumps their Priority from "low" to "ht" and
# refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
# timestamps to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 4-7 on each localized sheet correspond to the 4 files being re-handed-off.
$rows = 4..7

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column (G) is shared
    # between the languages - update it to the latest (de-de) generation time
    # so it stays a single shared value, same as before the edit.
    $overview.Range("G$r").Value = "2016-08-30 14:39:44"

    # zh-cn sheet: bump Priority (E) to "ht" and the Latest Handoff Datetime (H)
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-30 14:39:39"

    # de-de sheet: bump Priority (E) to "ht" and the Latest Handoff Datetime (H)
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-30 14:39:44"
}
